$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new data point for the "Smart Pointer to Implementation" column (L)
$ws.Range("L3").Value = 1046

# Move the active selection to L4, matching the saved cursor position
$ws.Range("L4").Select()

$wb.Save()
